$d = $word.ActiveDocument

# --- Hunk 1: merge the "A" / " mutation is a change..." runs into a single
#     run (text content is unchanged - only the run split goes away). ---
$found1 = $d.Content.Find.Execute(
    "A mutation is a change", $true, $false, $false, $false, $false,
    $true, 1, $false, "A mutation is a change", 2)

# Re-touch the whole (now merged) paragraph text so the surviving run keeps
# an (empty) <w:rPr/> element, matching the original authoring style.
$p5 = $d.Paragraphs(5).Range
$p5full = $d.Range($p5.Start, $p5.End - 1)
$p5full.Bold = 1
$p5full.Bold = 0

# --- Hunk 2: "ascii" -> "ansi" within "...saved as ascii non-unicode...",
#     but the resulting paragraph keeps the word split into three runs:
#     "...saved as a" | "nsi" | " non-unicode text files." ---
$rng = $d.Content
$found2 = $rng.Find.Execute(
    "ascii non-unicode", $true, $false, $false, $false, $false,
    $true, 1, $false, "ansi non-unicode", 2)

$p14 = $d.Paragraphs(14).Range
$pStart = $p14.Start
$pEnd = $p14.End - 1

$sub = $d.Range($pStart, $pEnd)
$subFound = $sub.Find.Execute(
    "nsi", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)

$seg1 = $d.Range($pStart, $sub.Start)
$seg2 = $d.Range($sub.Start, $sub.End)
$seg3 = $d.Range($sub.End, $pEnd)

# Toggling a character property and back forces Word to keep these three
# ranges as distinct runs (each with an empty <w:rPr/>) instead of Word
# re-coalescing them into one run.
$seg1.Bold = 1
$seg1.Bold = 0
$seg2.Bold = 1
$seg2.Bold = 0
$seg3.Bold = 1
$seg3.Bold = 0

Write-Host "Hunk1 found:" $found1 " Hunk2 found:" $found2 " nsi found:" $subFound
